# Rename the embedded logo pictures in the document's header/footer parts.
#
# The BTEC logo (header1.xml / header2.xml -> Word "first page" / "primary"
# headers) is currently named "image1.jpg" and must become "image2.jpg".
#
# The Pearson logo (footer1.xml / footer2.xml -> Word "first page" /
# "primary" footers) is currently named "image2.png" and must become
# "image1.png".
#
# Both the "first page" header/footer (wdHeaderFooterFirstPage, index 2)
# and the "primary"/default header/footer (wdHeaderFooterPrimary, index 1)
# carry an identical copy of the logo, so every section's Headers(1),
# Headers(2), Footers(1) and Footers(2) collection needs the same rename
# applied to its InlineShape.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # --- Headers: BTec logo, image1.jpg -> image2.jpg -------------------
    for ($i = 1; $i -le 2; $i++) {
        $hdr = $sec.Headers($i)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.Type -eq 3) {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    # --- Footers: Pearson logo, image2.png -> image1.png ----------------
    for ($i = 1; $i -le 2; $i++) {
        $ftr = $sec.Footers($i)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.Type -eq 3) {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
